$p = $ppt.ActivePresentation

# The "Development Process" slide is the last slide (slide 10); its body
# Content Placeholder (shape 2) gets its bullet text rewritten.
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# --- Paragraph 1: process arrow line -------------------------------------
$tr.Paragraphs(1).Text = "Requirement analysis -> Program Design -> Implementation -> AI integration ->Function improvements -> Verify"

# --- Paragraph 2: "Our project idea ..." -> "Idea from project2 and project3"
$tr.Paragraphs(2).Text = "Idea from project2 and project3"

# --- Paragraph 3: "Based on the second job ..." -> "Design a travel planner..."
$tr.Paragraphs(3).Text = "Design a travel planner closer to real-life travel. "

# Insert new paragraph right after paragraph 3: "Extended anytime search "
$null = $tr.Paragraphs(3).InsertAfter([char]13 + "Extended anytime search ")

# --- Paragraph 5 (was paragraph 4 "ChatGPT4.0 was used ..."): -> "Implement daily planner"
$tr.Paragraphs(5).Text = "Implement daily planner"

# Insert new paragraph right after paragraph 5: "ChatGPT4.0 helps neural network design and training "
$null = $tr.Paragraphs(5).InsertAfter([char]13 + "ChatGPT4.0 helps neural network design and training ")

# --- Paragraph 7 (was paragraph 5 "During the development ..."): -> "We are more familiar..."
$tr.Paragraphs(7).Text = "We are more familiar with search algorithms and neural network training"

# Insert new paragraph right after paragraph 7: "In future, we may extend ..."
$null = $tr.Paragraphs(7).InsertAfter([char]13 + "In future, we may extend the front-end of this project to implement a graphical interface")
